$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20, shifting rows 20:25 down to 21:26.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44518
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112022
$ws.Range("G20").Value = "Arveja Verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 620
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
